{"js": "// Remove the trailing \"Ver no Jupiter...\" / \"\u00a9 2020 ...\" footer block (and\n// the blank paragraph right before it) that the Jekyll build script used to\n// append after the bibliography entry \"Rio de Janeiro: Elsevier Editora, 2007.\".\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nconst markers = [\n  \"Ver no Jupiter Salvar em pdf Salvar em docx\",\n  \"\u00a9 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution\"\n];\n\n// Find the index of the first marker paragraph; the blank paragraph that\n// immediately precedes it (the spacer before the footer block) is removed\n// together with both marker paragraphs.\nlet firstMarkerIndex = -1;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text === markers[0]) {\n    firstMarkerIndex = i;\n    break;\n  }\n}\n\nconst toDelete = [];\nif (firstMarkerIndex !== -1) {\n  if (firstMarkerIndex - 1 >= 0 && paragraphs.items[firstMarkerIndex - 1].text === \"\") {\n    toDelete.push(paragraphs.items[firstMarkerIndex - 1]);\n  }\n  toDelete.push(paragraphs.items[firstMarkerIndex]);\n  if (firstMarkerIndex + 1 < paragraphs.items.length && paragraphs.items[firstMarkerIndex + 1].text === markers[1]) {\n    toDelete.push(paragraphs.items[firstMarkerIndex + 1]);\n  }\n}\n\nfor (const p of toDelete) {\n  p.delete();\n}\nawait context.sync();\n", "ps1": "# Remove the trailing \"Ver no Jupiter...\" / \"(c) 2020 ...\" footer block (and\n# the blank spacer paragraph right before it) that the Jekyll build script\n# used to append after the bibliography entry\n# \"Rio de Janeiro: Elsevier Editora, 2007.\".\n\n$d = $word.ActiveDocument\n\n$marker1 = \"Ver no Jupiter Salvar em pdf Salvar em docx\"\n$marker2 = [char]0x00A9 + \" 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution\"\n\n# Locate the first marker paragraph by its text (paragraph mark stripped).\n$idx1 = -1\n$count = $d.Paragraphs.Count\nfor ($i = 1; $i -le $count; $i++) {\n    $t = $d.Paragraphs($i).Range.Text.TrimEnd([char]13)\n    if ($t -eq $marker1) {\n        $idx1 = $i\n        break\n    }\n}\n\nif ($idx1 -gt 0) {\n    # Delete from the bottom up so earlier indices stay valid.\n    $idx2 = $idx1 + 1\n    if ($idx2 -le $d.Paragraphs.Count) {\n        $t2 = $d.Paragraphs($idx2).Range.Text.TrimEnd([char]13)\n        if ($t2 -eq $marker2) {\n            $d.Paragraphs($idx2).Range.Delete()\n        }\n    }\n\n    $d.Paragraphs($idx1).Range.Delete()\n\n    $idxBlank = $idx1 - 1\n    if ($idxBlank -ge 1) {\n        $tBlank = $d.Paragraphs($idxBlank).Range.Text.TrimEnd([char]13)\n        if ($tBlank -eq \"\") {\n            $d.Paragraphs($idxBlank).Range.Delete()\n        }\n    }\n}\n"}
